$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1": insert new order as row 14, push old "Tổng" row to row 15 ---
$wsPP1 = $wb.Worksheets.Item("Đơn phụ phẫu 1")

# Shift the old totals row (14) down to row 15 by writing it fresh at row 15,
# then overwrite row 14 with the new order data.
$wsPP1.Cells.Item(15, 1).Value = "Tổng"
$wsPP1.Cells.Item(15, 2).Value = 13
$wsPP1.Cells.Item(15, 3).Value = $null
$wsPP1.Cells.Item(15, 4).Value = $null
$wsPP1.Cells.Item(15, 5).Value = $null
$wsPP1.Cells.Item(15, 6).Value = $null
$wsPP1.Cells.Item(15, 7).Value = $null
$wsPP1.Cells.Item(15, 8).Value = $null
$wsPP1.Cells.Item(15, 9).Value = 900000

$wsPP1.Cells.Item(14, 1).Value = "HD-LUXURY"
$wsPP1.Cells.Item(14, 2).Value = 684
$wsPP1.Cells.Item(14, 3).NumberFormat = "@"
$wsPP1.Cells.Item(14, 3).Value = "08-25-2024"
$wsPP1.Cells.Item(14, 3).Style = "Normal"
$wsPP1.Cells.Item(14, 4).Value = "CẦN THƠ"
$wsPP1.Cells.Item(14, 5).Value = "Nguyễn Thị Như Ý"
$wsPP1.Cells.Item(14, 6).Value = "Cá nhân"
$wsPP1.Cells.Item(14, 7).Value = "Nâng mũi"
$wsPP1.Cells.Item(14, 8).Value = "Lâm Hoàng Phú"
$wsPP1.Cells.Item(14, 9).Value = 100000

# --- Sheet "Đơn thu nợ": update Chiết khấu sale chính (S column) values ---
$wsTN = $wb.Worksheets.Item("Đơn thu nợ")
$wsTN.Cells.Item(2, 19).Value = 90000
$wsTN.Cells.Item(3, 19).Value = 60000.00000000001
$wsTN.Cells.Item(4, 19).Value = 150000

# --- Sheet "Lương": update computed payroll summary values ---
$wsLuong = $wb.Worksheets.Item("Lương")
$wsLuong.Cells.Item(2, 2).Value = 22
$wsLuong.Cells.Item(3, 2).Value = 2553571.428571429
$wsLuong.Cells.Item(8, 2).Value = 900000
$wsLuong.Cells.Item(10, 2).Value = 150000
$wsLuong.Cells.Item(34, 2).Value = 803571.4285714286
$wsLuong.Cells.Item(37, 2).Value = 903571.4285714286
